# ---------------------------------------------------------------------------
# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (fund-holding detail) right before the
# "总计" (totals) summary sheet, populates it with the Q1-2022 holding detail
# rows, and prepends a matching summary row to the "总计" sheet.
#
# NOTE: worksheet/range handles in this host are position-based, so any handle
# captured *before* a sheet is inserted/removed can silently start pointing at
# a different sheet once positions shift. To stay safe we always look sheets
# back up **by name** right before we use them, instead of reusing a handle
# across an Add()/Insert() call.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# A sheet with the same 基金代码/基金名称/... layout already exists (the most
# recent quarter, "2021-Q4"); reuse it as a style template so the new sheet
# matches the existing look (bold/bordered header row + bold index column).
$templateSheet = $wb.Worksheets.Item("2021-Q4")

# "总计" is currently the last sheet; insert "2022-Q1" right before it.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($lastSheet)
$ws.Name = "2022-Q1"

# Copy header/index-column formatting from the template sheet.
$templateSheet.Range("B1:H1").Copy() | Out-Null
$ws.Range("B1:H1").PasteSpecial(-4122)
$templateSheet.Range("A2:A11").Copy() | Out-Null
$ws.Range("A2:A17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- 2. Header row ----------------------------------------------------------
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# --- 3. Fund holding detail rows --------------------------------------------
$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "'011304"
$ws.Cells.Item(2, 3).Value = "工银瑞信创新成长混合A"
$ws.Cells.Item(2, 4).Value = "'44.47"
$ws.Cells.Item(2, 5).Value = "'83.84"
$ws.Cells.Item(2, 6).Value = "'3.26"
$ws.Cells.Item(2, 7).Value = "'1.4497"
$ws.Cells.Item(2, 8).Value = 4

$ws.Cells.Item(3, 1).Value = 1
$ws.Cells.Item(3, 2).Value = "'166005"
$ws.Cells.Item(3, 3).Value = "中欧价值发现混合 -A"
$ws.Cells.Item(3, 4).Value = "'43.52"
$ws.Cells.Item(3, 5).Value = "'93.97"
$ws.Cells.Item(3, 6).Value = "'3.28"
$ws.Cells.Item(3, 7).Value = "'1.4275"
$ws.Cells.Item(3, 8).Value = 7

$ws.Cells.Item(4, 1).Value = 2
$ws.Cells.Item(4, 2).Value = "'001882"
$ws.Cells.Item(4, 3).Value = "中欧价值发现混合 -E"
$ws.Cells.Item(4, 4).Value = "'43.52"
$ws.Cells.Item(4, 5).Value = "'93.97"
$ws.Cells.Item(4, 6).Value = "'3.28"
$ws.Cells.Item(4, 7).Value = "'1.4275"
$ws.Cells.Item(4, 8).Value = 7

$ws.Cells.Item(5, 1).Value = 3
$ws.Cells.Item(5, 2).Value = "'001409"
$ws.Cells.Item(5, 3).Value = "工银瑞信互联网加股票"
$ws.Cells.Item(5, 4).Value = "'34.72"
$ws.Cells.Item(5, 5).Value = "'82.35"
$ws.Cells.Item(5, 6).Value = "'2.76"
$ws.Cells.Item(5, 7).Value = "'0.9583"
$ws.Cells.Item(5, 8).Value = 9

$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "'001810"
$ws.Cells.Item(6, 3).Value = "中欧潜力价值灵活配置混合A"
$ws.Cells.Item(6, 4).Value = "'28.67"
$ws.Cells.Item(6, 5).Value = "'94.05"
$ws.Cells.Item(6, 6).Value = "'3.33"
$ws.Cells.Item(6, 7).Value = "'0.9547"
$ws.Cells.Item(6, 8).Value = 7

$ws.Cells.Item(7, 1).Value = 5
$ws.Cells.Item(7, 2).Value = "'000263"
$ws.Cells.Item(7, 3).Value = "工银瑞信信息产业混合A"
$ws.Cells.Item(7, 4).Value = "'27.45"
$ws.Cells.Item(7, 5).Value = "'83.76"
$ws.Cells.Item(7, 6).Value = "'2.94"
$ws.Cells.Item(7, 7).Value = "'0.8070"
$ws.Cells.Item(7, 8).Value = 6

$ws.Cells.Item(8, 1).Value = 6
$ws.Cells.Item(8, 2).Value = "'004232"
$ws.Cells.Item(8, 3).Value = "中欧价值发现混合 -C"
$ws.Cells.Item(8, 4).Value = "'10.98"
$ws.Cells.Item(8, 5).Value = "'93.97"
$ws.Cells.Item(8, 6).Value = "'3.28"
$ws.Cells.Item(8, 7).Value = "'0.3601"
$ws.Cells.Item(8, 8).Value = 7

$ws.Cells.Item(9, 1).Value = 7
$ws.Cells.Item(9, 2).Value = "'005764"
$ws.Cells.Item(9, 3).Value = "中欧潜力价值灵活配置混合C"
$ws.Cells.Item(9, 4).Value = "'3.43"
$ws.Cells.Item(9, 5).Value = "'94.05"
$ws.Cells.Item(9, 6).Value = "'3.33"
$ws.Cells.Item(9, 7).Value = "'0.1142"
$ws.Cells.Item(9, 8).Value = 7

$ws.Cells.Item(10, 1).Value = 8
$ws.Cells.Item(10, 2).Value = "'001891"
$ws.Cells.Item(10, 3).Value = "中欧成长优选回报灵活配置混合E"
$ws.Cells.Item(10, 4).Value = "'2.97"
$ws.Cells.Item(10, 5).Value = "'94.42"
$ws.Cells.Item(10, 6).Value = "'3.30"
$ws.Cells.Item(10, 7).Value = "'0.0980"
$ws.Cells.Item(10, 8).Value = 6

$ws.Cells.Item(11, 1).Value = 9
$ws.Cells.Item(11, 2).Value = "'166020"
$ws.Cells.Item(11, 3).Value = "中欧成长优选回报灵活配置混合A"
$ws.Cells.Item(11, 4).Value = "'2.97"
$ws.Cells.Item(11, 5).Value = "'94.42"
$ws.Cells.Item(11, 6).Value = "'3.30"
$ws.Cells.Item(11, 7).Value = "'0.0980"
$ws.Cells.Item(11, 8).Value = 6

$ws.Cells.Item(12, 1).Value = 10
$ws.Cells.Item(12, 2).Value = "'011474"
$ws.Cells.Item(12, 3).Value = "工银瑞信信息产业混合C"
$ws.Cells.Item(12, 4).Value = "'2.67"
$ws.Cells.Item(12, 5).Value = "'83.76"
$ws.Cells.Item(12, 6).Value = "'2.94"
$ws.Cells.Item(12, 7).Value = "'0.0785"
$ws.Cells.Item(12, 8).Value = 6

$ws.Cells.Item(13, 1).Value = 11
$ws.Cells.Item(13, 2).Value = "'002581"
$ws.Cells.Item(13, 3).Value = "招商丰凯灵活配置混合A"
$ws.Cells.Item(13, 4).Value = "'5.76"
$ws.Cells.Item(13, 5).Value = "'40.60"
$ws.Cells.Item(13, 6).Value = "'1.20"
$ws.Cells.Item(13, 7).Value = "'0.0691"
$ws.Cells.Item(13, 8).Value = 9

$ws.Cells.Item(14, 1).Value = 12
$ws.Cells.Item(14, 2).Value = "'011305"
$ws.Cells.Item(14, 3).Value = "工银瑞信创新成长混合C"
$ws.Cells.Item(14, 4).Value = "'1.61"
$ws.Cells.Item(14, 5).Value = "'83.84"
$ws.Cells.Item(14, 6).Value = "'3.26"
$ws.Cells.Item(14, 7).Value = "'0.0525"
$ws.Cells.Item(14, 8).Value = 4

$ws.Cells.Item(15, 1).Value = 13
$ws.Cells.Item(15, 2).Value = "'487016"
$ws.Cells.Item(15, 3).Value = "工银瑞信灵活配置混合A"
$ws.Cells.Item(15, 4).Value = "'2.35"
$ws.Cells.Item(15, 5).Value = "'73.70"
$ws.Cells.Item(15, 6).Value = "'1.66"
$ws.Cells.Item(15, 7).Value = "'0.0390"
$ws.Cells.Item(15, 8).Value = 7

$ws.Cells.Item(16, 1).Value = 14
$ws.Cells.Item(16, 2).Value = "'002582"
$ws.Cells.Item(16, 3).Value = "招商丰凯灵活配置混合C"
$ws.Cells.Item(16, 4).Value = "'1.46"
$ws.Cells.Item(16, 5).Value = "'40.60"
$ws.Cells.Item(16, 6).Value = "'1.20"
$ws.Cells.Item(16, 7).Value = "'0.0175"
$ws.Cells.Item(16, 8).Value = 9

$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = "'001428"
$ws.Cells.Item(17, 3).Value = "工银瑞信灵活配置混合B"
$ws.Cells.Item(17, 4).Value = "'0.00"
$ws.Cells.Item(17, 5).Value = "'73.70"
$ws.Cells.Item(17, 6).Value = "'1.66"
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 7

# --- 4. Prepend a "2022-Q1" summary row to the "总计" sheet -----------------
# Re-fetch "总计" by name: the worksheet that used to be $lastSheet has now
# shifted from position 6 to position 7 because of the insertion above.
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()
$totalSheet.Range("B2:D2").ClearFormats()
$totalSheet.Range("A3").Copy() | Out-Null
$totalSheet.Range("A2").PasteSpecial(-4122)
$excel.CutCopyMode = 0
$totalSheet.Cells.Item(2, 2).Value = "2022-Q1"
$totalSheet.Cells.Item(2, 3).Value = 16
$totalSheet.Cells.Item(2, 4).Value = 7.95

# The numeric index in column A is a plain running counter (0, 1, 2, ...),
# not a formula, so re-number every data row from scratch top to bottom.
$totalSheet.Cells.Item(2, 1).Value = 0
$totalSheet.Cells.Item(3, 1).Value = 1
$totalSheet.Cells.Item(4, 1).Value = 2
$totalSheet.Cells.Item(5, 1).Value = 3
$totalSheet.Cells.Item(6, 1).Value = 4
$totalSheet.Cells.Item(7, 1).Value = 5

